# Making sure the data input is equal to the base case
#
# Fuel_Cost_Absolute: Wind (B5) and PV (B6) cost drop from 5 -> 0, to match
# the base-case data input. This ripples into Python_Gen_N_Data!D7:D17
# (INDEX/MATCH lookups) automatically on recalculation.

$wb = $excel.ActiveWorkbook

$fuelCost = $wb.Worksheets.Item("Fuel_Cost_Absolute")
$fuelCost.Range("B5").Value = 0
$fuelCost.Range("B6").Value = 0

# Update the remembered selection on that sheet.
$fuelCost.Range("B2").Select()

# The active tab moves on to Python_Gen_N_Data, with a new remembered
# selection there too.
$pyGenData = $wb.Worksheets.Item("Python_Gen_N_Data")
$pyGenData.Activate()
$pyGenData.Range("H19").Select()
